$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Icam1"
$ws.Cells.Item(2,3).Value = "Itgal"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 14.452709
$ws.Cells.Item(2,8).Value = 43.358127
$ws.Cells.Item(2,9).Value = 0.1476906377370901
$ws.Cells.Item(2,10).Value = 0.1476906377370901
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.9976189999999999
$ws.Cells.Item(2,14).Value = 2.992857
$ws.Cells.Item(2,15).Value = 0.06139866523632875
$ws.Cells.Item(2,16).Value = 0.06139866523632876
$ws.Cells.Item(2,17).Value = 14.418297099871
$ws.Cells.Item(2,18).Value = 129.764673898839
$ws.Cells.Item(2,19).Value = 0.009068008024959497
$ws.Cells.Item(2,20).Value = 0.0090680080249595

# Row 3: ECs -> M2
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Icam1"
$ws.Cells.Item(3,3).Value = "Itgal"
$ws.Cells.Item(3,4).Value = "M2"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 14.452709
$ws.Cells.Item(3,8).Value = 43.358127
$ws.Cells.Item(3,9).Value = 0.1476906377370901
$ws.Cells.Item(3,10).Value = 0.1476906377370901
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 14.969586
$ws.Cells.Item(3,14).Value = 44.908758
$ws.Cells.Item(3,15).Value = 0.9213062296732858
$ws.Cells.Item(3,16).Value = 0.9213062296732859
$ws.Cells.Item(3,17).Value = 216.351070308474
$ws.Cells.Item(3,18).Value = 1947.159632776266
$ws.Cells.Item(3,19).Value = 0.1360683046116016
$ws.Cells.Item(3,20).Value = 0.1360683046116016

# Row 4: ECs -> sCs
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Icam1"
$ws.Cells.Item(4,3).Value = "Itgal"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 14.452709
$ws.Cells.Item(4,8).Value = 43.358127
$ws.Cells.Item(4,9).Value = 0.1476906377370901
$ws.Cells.Item(4,10).Value = 0.1476906377370901
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.2810146666666667
$ws.Cells.Item(4,14).Value = 0.843044
$ws.Cells.Item(4,15).Value = 0.01729510509038539
$ws.Cells.Item(4,16).Value = 0.01729510509038539
$ws.Cells.Item(4,17).Value = 4.061423202065333
$ws.Cells.Item(4,18).Value = 36.552808818588
$ws.Cells.Item(4,19).Value = 0.002554325100529012
$ws.Cells.Item(4,20).Value = 0.002554325100529012

# Row 5: FAPs -> FAPs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Icam1"
$ws.Cells.Item(5,3).Value = "Itgal"
$ws.Cells.Item(5,4).Value = "FAPs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 45.91529066666667
$ws.Cells.Item(5,8).Value = 137.745872
$ws.Cells.Item(5,9).Value = 0.4692032864180593
$ws.Cells.Item(5,10).Value = 0.4692032864180593
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.9976189999999999
$ws.Cells.Item(5,14).Value = 2.992857
$ws.Cells.Item(5,15).Value = 0.06139866523632875
$ws.Cells.Item(5,16).Value = 0.06139866523632876
$ws.Cells.Item(5,17).Value = 45.80596635958933
$ws.Cells.Item(5,18).Value = 412.253697236304
$ws.Cells.Item(5,19).Value = 0.0288084555105677
$ws.Cells.Item(5,20).Value = 0.0288084555105677

# Row 6: FAPs -> M2
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Icam1"
$ws.Cells.Item(6,3).Value = "Itgal"
$ws.Cells.Item(6,4).Value = "M2"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 45.91529066666667
$ws.Cells.Item(6,8).Value = 137.745872
$ws.Cells.Item(6,9).Value = 0.4692032864180593
$ws.Cells.Item(6,10).Value = 0.4692032864180593
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 14.969586
$ws.Cells.Item(6,14).Value = 44.908758
$ws.Cells.Item(6,15).Value = 0.9213062296732858
$ws.Cells.Item(6,16).Value = 0.9213062296732859
$ws.Cells.Item(6,17).Value = 687.332892349664
$ws.Cells.Item(6,18).Value = 6185.996031146977
$ws.Cells.Item(6,19).Value = 0.432279910760137
$ws.Cells.Item(6,20).Value = 0.4322799107601371

# Row 7: FAPs -> sCs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Icam1"
$ws.Cells.Item(7,3).Value = "Itgal"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 45.91529066666667
$ws.Cells.Item(7,8).Value = 137.745872
$ws.Cells.Item(7,9).Value = 0.4692032864180593
$ws.Cells.Item(7,10).Value = 0.4692032864180593
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 0.2810146666666667
$ws.Cells.Item(7,14).Value = 0.843044
$ws.Cells.Item(7,15).Value = 0.01729510509038539
$ws.Cells.Item(7,16).Value = 0.01729510509038539
$ws.Cells.Item(7,17).Value = 12.90287010159645
$ws.Cells.Item(7,18).Value = 116.125830914368
$ws.Cells.Item(7,19).Value = 0.008114920147354531
$ws.Cells.Item(7,20).Value = 0.008114920147354533

# Row 8: M2 -> FAPs
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Icam1"
$ws.Cells.Item(8,3).Value = "Itgal"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 36.015061
$ws.Cells.Item(8,8).Value = 108.045183
$ws.Cells.Item(8,9).Value = 0.3680339324088102
$ws.Cells.Item(8,10).Value = 0.3680339324088103
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.9976189999999999
$ws.Cells.Item(8,14).Value = 2.992857
$ws.Cells.Item(8,15).Value = 0.06139866523632875
$ws.Cells.Item(8,16).Value = 0.06139866523632876
$ws.Cells.Item(8,17).Value = 35.929309139759
$ws.Cells.Item(8,18).Value = 323.363782257831
$ws.Cells.Item(8,19).Value = 0.02259679221157818
$ws.Cells.Item(8,20).Value = 0.02259679221157819

# Row 9: M2 -> M2
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Icam1"
$ws.Cells.Item(9,3).Value = "Itgal"
$ws.Cells.Item(9,4).Value = "M2"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 36.015061
$ws.Cells.Item(9,8).Value = 108.045183
$ws.Cells.Item(9,9).Value = 0.3680339324088102
$ws.Cells.Item(9,10).Value = 0.3680339324088103
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 14.969586
$ws.Cells.Item(9,14).Value = 44.908758
$ws.Cells.Item(9,15).Value = 0.9213062296732858
$ws.Cells.Item(9,16).Value = 0.9213062296732859
$ws.Cells.Item(9,17).Value = 539.130552934746
$ws.Cells.Item(9,18).Value = 4852.174976412714
$ws.Cells.Item(9,19).Value = 0.3390719546593938
$ws.Cells.Item(9,20).Value = 0.3390719546593939

# Row 10: M2 -> sCs
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Icam1"
$ws.Cells.Item(10,3).Value = "Itgal"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 36.015061
$ws.Cells.Item(10,8).Value = 108.045183
$ws.Cells.Item(10,9).Value = 0.3680339324088102
$ws.Cells.Item(10,10).Value = 0.3680339324088103
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.2810146666666667
$ws.Cells.Item(10,14).Value = 0.843044
$ws.Cells.Item(10,15).Value = 0.01729510509038539
$ws.Cells.Item(10,16).Value = 0.01729510509038539
$ws.Cells.Item(10,17).Value = 10.12076036189467
$ws.Cells.Item(10,18).Value = 91.08684325705201
$ws.Cells.Item(10,19).Value = 0.006365185537838167
$ws.Cells.Item(10,20).Value = 0.006365185537838168

# Row 11: sCs -> FAPs
$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Icam1"
$ws.Cells.Item(11,3).Value = "Itgal"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 1.474929666666667
$ws.Cells.Item(11,8).Value = 4.424789000000001
$ws.Cells.Item(11,9).Value = 0.01507214343604052
$ws.Cells.Item(11,10).Value = 0.01507214343604052
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.9976189999999999
$ws.Cells.Item(11,14).Value = 2.992857
$ws.Cells.Item(11,15).Value = 0.06139866523632875
$ws.Cells.Item(11,16).Value = 0.06139866523632876
$ws.Cells.Item(11,17).Value = 1.471417859130334
$ws.Cells.Item(11,18).Value = 13.242760732173
$ws.Cells.Item(11,19).Value = 0.0009254094892233818
$ws.Cells.Item(11,20).Value = 0.0009254094892233819

# Row 12: sCs -> M2
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Icam1"
$ws.Cells.Item(12,3).Value = "Itgal"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 1.474929666666667
$ws.Cells.Item(12,8).Value = 4.424789000000001
$ws.Cells.Item(12,9).Value = 0.01507214343604052
$ws.Cells.Item(12,10).Value = 0.01507214343604052
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 14.969586
$ws.Cells.Item(12,14).Value = 44.908758
$ws.Cells.Item(12,15).Value = 0.9213062296732858
$ws.Cells.Item(12,16).Value = 0.9213062296732859
$ws.Cells.Item(12,17).Value = 22.079086489118
$ws.Cells.Item(12,18).Value = 198.711778402062
$ws.Cells.Item(12,19).Value = 0.01388605964215346
$ws.Cells.Item(12,20).Value = 0.01388605964215346

# Row 13: sCs -> sCs
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Icam1"
$ws.Cells.Item(13,3).Value = "Itgal"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 1.474929666666667
$ws.Cells.Item(13,8).Value = 4.424789000000001
$ws.Cells.Item(13,9).Value = 0.01507214343604052
$ws.Cells.Item(13,10).Value = 0.01507214343604052
$ws.Cells.Item(13,11).Value = 2
$ws.Cells.Item(13,12).Value = 0.6666666666666666
$ws.Cells.Item(13,13).Value = 0.2810146666666667
$ws.Cells.Item(13,14).Value = 0.843044
$ws.Cells.Item(13,15).Value = 0.01729510509038539
$ws.Cells.Item(13,16).Value = 0.01729510509038539
$ws.Cells.Item(13,17).Value = 0.4144768686351112
$ws.Cells.Item(13,18).Value = 3.730291817716
$ws.Cells.Item(13,19).Value = 0.0002606743046636832
$ws.Cells.Item(13,20).Value = 0.0002606743046636832
